# All Result ready to start write
# The sheet's leftmost column (old "A", an 11/15 taxon-count helper column)
# is being dropped entirely: every other column shifts one place to the
# left (B->A, C->B, D->C, E->D, F->E) and the dimension shrinks from
# A1:F3 to A1:E3.
#
# Deleting the entire column A accomplishes exactly this shift (values,
# shared-string references and styles all move left together), which is
# the same result produced by the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:A3").EntireColumn.Delete()
